$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.418.44"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = "'2.376.85"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +5.67%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'235.89"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('D6').Value = "'0.648"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('D7').Value = "'71.19"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +12.98%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('D11').Value = "'56.99"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = "'26.51"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').Value = "'2.724.15"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.49%  '
$ws.Range('D14').Value = "'0.106"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').Value = "'15.78"
$ws.Range('D15').Style = "Normal"
$ws.Range('E16').Value = '  +2.92%  '
$ws.Range('E17').Value = '  +3.56%  '
$ws.Range('D18').Value = "'2.377.16"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +5.97%  '
$ws.Range('D19').Value = "'43.425.91"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = "'0.0₃0990"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').Value = "'6.36"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.94%  '
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').Value = "'251.62"
$ws.Range('D23').Style = "Normal"
$ws.Range('E24').Value = '  +19.12%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = "'2.47"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('D27').Value = "'23.06"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +9.95%  '
$ws.Range('D28').Value = "'10.02"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').Value = "'2.23"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').Value = "'174.31"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = "'1.54"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.15%  '
$ws.Range('E32').Value = '  -8.47%  '
$ws.Range('D33').Value = "'0.127"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.12%  '
$ws.Range('D34').Value = "'5.00"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.61%  '
$ws.Range('D35').Value = "'0.0692"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Value = "'6.61"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = "'2.46"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +8.20%  '
$ws.Range('D39').Value = "'3.65"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = "'8.97"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.97%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').Value = "'18.59"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +9.23%  '
$ws.Range('E44').Value = '  +10.45%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'99.87"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.58%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = "'4.52"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.40%  '
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').Value = "'0.0946"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').Value = "'1.453.67"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').Value = "'2.599.25"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.78%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').Value = "'0.000201"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -9.45%  '
